$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5341.0356
$ws.Range("I40").Value = 4106.04
$ws.Range("J40").Value = 15632.667
$ws.Range("K40").Value = 4106.04
$ws.Range("L40").Value = 15632.667
$ws.Range("M40").Value = -3931.04
$ws.Range("N40").Value = -15982.667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 3152.8948
$ws.Range("I107").Value = 1913.9333
$ws.Range("K107").Value = 1913.9333
$ws.Range("M107").Value = 6.066700000000083

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 4762.385
$ws.Range("I115").Value = 4991.909
$ws.Range("K115").Value = 14975.727
$ws.Range("M115").Value = -13408.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6653.1113
$ws.Range("I116").Value = 4978.6665
$ws.Range("J116").Value = 10002
$ws.Range("K116").Value = 4978.6665
$ws.Range("L116").Value = 10002
$ws.Range("M116").Value = -1536.6665
$ws.Range("N116").Value = -16886

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 9838.233
$ws.Range("I132").Value = 4766.3335
$ws.Range("J132").Value = 30125.834
$ws.Range("K132").Value = 14299.0005
$ws.Range("L132").Value = 90377.50199999999
$ws.Range("M132").Value = -11769.0005
$ws.Range("N132").Value = -95437.50199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4253.56
$ws.Range("I32").Value = 3125.9307
$ws.Range("K32").Value = 3125.9307
$ws.Range("M32").Value = -2838.9307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 15197.8
$ws.Range("I45").Value = 27481.6
$ws.Range("K45").Value = 27481.6
$ws.Range("M45").Value = -27104.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2782.348
$ws.Range("I74").Value = 2845.2273
$ws.Range("K74").Value = 2845.2273
$ws.Range("M74").Value = -1971.2273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2782.348
$ws.Range("I77").Value = 2845.2273
$ws.Range("K77").Value = 14226.1365
$ws.Range("M77").Value = -9858.136500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 45410.695
$ws.Range("I97").Value = 2394.3333
$ws.Range("J97").Value = 200269.6
$ws.Range("K97").Value = 2394.3333
$ws.Range("L97").Value = 200269.6
$ws.Range("M97").Value = -1898.3333
$ws.Range("N97").Value = -201261.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 56124.11
$ws.Range("I22").Value = 639.625
$ws.Range("K22").Value = 639.625
$ws.Range("M22").Value = -466.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1926.9
$ws.Range("I94").Value = 1659.0834
$ws.Range("K94").Value = 1659.0834
$ws.Range("M94").Value = -1208.0834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2762.5
$ws.Range("I107").Value = 2894.2856
$ws.Range("K107").Value = 2894.2856
$ws.Range("M107").Value = -974.2856000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6693.7036
$ws.Range("I134").Value = 6629.24
$ws.Range("K134").Value = 19887.72
$ws.Range("M134").Value = -17352.72

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1669.5264
$ws.Range("I31").Value = 1045.0625
$ws.Range("K31").Value = 1045.0625
$ws.Range("M31").Value = -750.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1669.5264
$ws.Range("I34").Value = 1045.0625
$ws.Range("K34").Value = 1045.0625
$ws.Range("M34").Value = -843.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8272.280000000001
$ws.Range("I58").Value = 7183.25
$ws.Range("J58").Value = 10208.333
$ws.Range("K58").Value = 7183.25
$ws.Range("L58").Value = 10208.333
$ws.Range("M58").Value = -6980.25
$ws.Range("N58").Value = -10614.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2441.2
$ws.Range("J94").Value = 2617.5454
$ws.Range("L94").Value = 2617.5454
$ws.Range("N94").Value = -3519.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7366.647
$ws.Range("I99").Value = 5128.5
$ws.Range("J99").Value = 9356.111000000001
$ws.Range("K99").Value = 5128.5
$ws.Range("L99").Value = 9356.111000000001
$ws.Range("M99").Value = -3630.5
$ws.Range("N99").Value = -12352.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1484.9231
$ws.Range("J107").Value = 1578.7
$ws.Range("L107").Value = 1578.7
$ws.Range("N107").Value = -5418.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3133.3044
$ws.Range("J122").Value = 3388.4443
$ws.Range("L122").Value = 10165.3329
$ws.Range("N122").Value = -15065.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 7366.647
$ws.Range("I126").Value = 5128.5
$ws.Range("J126").Value = 9356.111000000001
$ws.Range("K126").Value = 15385.5
$ws.Range("L126").Value = 28068.333
$ws.Range("M126").Value = -12915.5
$ws.Range("N126").Value = -33008.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9138.76
$ws.Range("I132").Value = 9117.15
$ws.Range("K132").Value = 27351.45
$ws.Range("M132").Value = -24821.45

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8272.280000000001
$ws.Range("I136").Value = 7183.25
$ws.Range("J136").Value = 10208.333
$ws.Range("K136").Value = 21549.75
$ws.Range("L136").Value = 30624.999
$ws.Range("M136").Value = -18999.75
$ws.Range("N136").Value = -35724.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1072.3846
$ws.Range("J5").Value = 734.3333
$ws.Range("L5").Value = 2202.9999
$ws.Range("N5").Value = -2426.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 326.25
$ws.Range("J33").Value = 429
$ws.Range("L33").Value = 2574
$ws.Range("N33").Value = -3140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 9766.416999999999
$ws.Range("I68").Value = 425.25
$ws.Range("J68").Value = 14437
$ws.Range("K68").Value = 1275.75
$ws.Range("L68").Value = 43311
$ws.Range("M68").Value = -464.75
$ws.Range("N68").Value = -44933

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 9766.416999999999
$ws.Range("I71").Value = 425.25
$ws.Range("J71").Value = 14437
$ws.Range("K71").Value = 3827.25
$ws.Range("L71").Value = 129933
$ws.Range("M71").Value = 228.75
$ws.Range("N71").Value = -138045

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2258.5
$ws.Range("J129").Value = 4256.1113
$ws.Range("L129").Value = 12768.3339
$ws.Range("N129").Value = -22768.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1072.3846
$ws.Range("J135").Value = 734.3333
$ws.Range("L135").Value = 6608.9997
$ws.Range("N135").Value = -11678.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 10937.875
$ws.Range("I102").Value = 13378.6
$ws.Range("J102").Value = 6870
$ws.Range("K102").Value = 13378.6
$ws.Range("L102").Value = 6870
$ws.Range("M102").Value = -11756.6
$ws.Range("N102").Value = -10114

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 767.5
$ws.Range("I107").Value = 842
$ws.Range("K107").Value = 842
$ws.Range("M107").Value = 1078

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5999
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6799.2666
$ws.Range("I122").Value = 6922.3076
$ws.Range("J122").Value = 5999.5
$ws.Range("K122").Value = 20766.9228
$ws.Range("L122").Value = 17998.5
$ws.Range("M122").Value = -18316.9228
$ws.Range("N122").Value = -22898.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5257
$ws.Range("I40").Value = 5257
$ws.Range("K40").Value = 5257
$ws.Range("M40").Value = -5121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 422.33
$ws.Range("I82").Value = 411.02063
$ws.Range("J82").Value = 788
$ws.Range("K82").Value = 411.02063
$ws.Range("L82").Value = 788
$ws.Range("M82").Value = -50.02062999999998
$ws.Range("N82").Value = -1510

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 422.33
$ws.Range("I85").Value = 411.02063
$ws.Range("J85").Value = 788
$ws.Range("K85").Value = 411.02063
$ws.Range("L85").Value = 788
$ws.Range("M85").Value = 836.97937
$ws.Range("N85").Value = -3284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 49999
$ws.Range("I42").Value = 49999
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 49999
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -49621
$ws.Range("N42").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 26100
$ws.Range("I62").Value = 15002
$ws.Range("K62").Value = 15002
$ws.Range("M62").Value = -14378

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 26100
$ws.Range("I65").Value = 15002
$ws.Range("K65").Value = 75010
$ws.Range("M65").Value = -71890

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2744.2856
$ws.Range("I96").Value = 3319
$ws.Range("K96").Value = 3319
$ws.Range("M96").Value = -1946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 774.5
$ws.Range("I113").Value = 442.4
$ws.Range("K113").Value = 1327.2
$ws.Range("M113").Value = 842.8000000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 8342365.5
$ws.Range("I126").Value = 11368572
$ws.Range("J126").Value = 20298.5
$ws.Range("K126").Value = 34105716
$ws.Range("L126").Value = 60895.5
$ws.Range("M126").Value = -34103246
$ws.Range("N126").Value = -65835.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1344.2295
$ws.Range("I132").Value = 1322.017
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 3966.051
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -1436.051
$ws.Range("N132").Value = -11058.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1608.2449
$ws.Range("I136").Value = 1609.375
$ws.Range("K136").Value = 4828.125
$ws.Range("M136").Value = -2278.125
